$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Datos de las sales volatiles a anadir (filas 8-13)
$datos = @(
    @("Siliciuro", "Si"),
    @("Carburo", "C "),
    @("Nitruro", "N"),
    @("Fosfuro", "P"),
    @("Arseniuro", "As"),
    @("Boruro", "B ")
)

$row = 8
foreach ($par in $datos) {
    $ws.Cells.Item($row, 1).Value = $par[0]
    $ws.Cells.Item($row, 2).Value = $par[1]
    $row = $row + 1
}

# Filas 11 y 12 (Fosfuro y Arseniuro) llevan ajuste de texto (wrap) en la columna A
$ws.Range("A11:A12").WrapText = $true
$ws.Range("A11:A12").Font.Name = "Arial"

# Dejar seleccionado el nuevo bloque de datos agregado
$ws.Range("A8:B13").Select()
